# Add new registrant rows (5-10) to the active sheet, mirroring the
# existing "Apellidos / Nombres / Colegio / Departamento / Provincia" table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("Claure Reyes",     "Milenia",      "Centro Cultural Anglo Americano Cochabamba", "Cochabamba", "Cercado"),
    @("Argote Gomez",     "Noelia",       "Colegio Particular Hispano Boliviano",       "Cochabamba", "Cercado"),
    @("Azogue Aranibar",  "Fernando",     "Colegio Particular Hispano Boliviano",       "Cochabamba", "Cercado"),
    @("Galvez Perez",     "Rodrigo",      "Centro Cultural Anglo Americano Cochabamba", "Cochabamba", "Sacaba"),
    @("Gonzales Medrano", "Emilio Jorge", "Colegio Loyola",                             "Cochabamba", "Cercado"),
    @("Vega Salinas",     "Juan Pablo",   "Colegio Albert Einstein",                    "Cochabamba", "Cercado")
)

$startRow = 5
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $rowIndex = $startRow + $i
    $rowData = $newRows[$i]
    for ($j = 0; $j -lt $rowData.Count; $j++) {
        $colIndex = $j + 1
        $ws.Cells.Item($rowIndex, $colIndex).Value = $rowData[$j]
    }
}
